$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.974.42'
$ws.Range('E2').Value = '  +1.42%  '
$ws.Range('D3').Value = '3.417.60'
$ws.Range('E3').Value = '  +0.71%  '
$ws.Range('E4').Value = '  -0.24%  '
$ws.Range('D5').Value = '406.28'
$ws.Range('E5').Value = '  -0.17%  '
$ws.Range('D6').Value = '132.09'
$ws.Range('E6').Value = '  +3.19%  '
$ws.Range('D7').Value = '0.593'
$ws.Range('E7').Value = '  -0.91%  '
$ws.Range('E8').Value = '  -0.32%  '
$ws.Range('E9').Value = '  +2.82%  '
$ws.Range('D10').Value = '0.132'
$ws.Range('E10').Value = '  +4.06%  '
$ws.Range('E11').Value = '  -0.68%  '
$ws.Range('E12').Value = '  -0.35%  '
$ws.Range('D13').Value = '19.99'
$ws.Range('E13').Value = '  +1.80%  '
$ws.Range('E14').Value = '  -1.36%  '
$ws.Range('D15').Value = '3.412.62'
$ws.Range('E15').Value = '  -0.44%  '
$ws.Range('D16').Value = '11.76'
$ws.Range('E16').Value = '  +1.09%  '
$ws.Range('D17').Value = '61.922.61'
$ws.Range('E17').Value = '  +0.83%  '
$ws.Range('E18').Value = '  +0.34%  '
$ws.Range('D19').Value = '0.0000148'
$ws.Range('E19').Value = '  +10.38%  '
$ws.Range('E20').Value = '  -1.91%  '
$ws.Range('D21').Value = '83.93'
$ws.Range('E21').Value = '  +1.53%  '
$ws.Range('D22').Value = '313.66'
$ws.Range('E22').Value = '  +2.36%  '
$ws.Range('D23').Value = '12.84'
$ws.Range('E23').Value = '  -0.93%  '
$ws.Range('D24').Value = '3.15'
$ws.Range('E24').Value = '  -0.27%  '
$ws.Range('D25').Value = '4.76'
$ws.Range('E25').Value = '  +3.42%  '
$ws.Range('E26').Value = '  +0.28%  '
$ws.Range('D27').Value = '8.00'
$ws.Range('E27').Value = '  +7.40%  '
$ws.Range('D28').Value = '8.13'
$ws.Range('E28').Value = '  -5.62%  '
$ws.Range('D29').Value = '2.76'
$ws.Range('E29').Value = '  +8.10%  '
$ws.Range('D30').Value = '0.173'
$ws.Range('E30').Value = '  +0.11%  '
$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D31').Value = '43.70'
$ws.Range('E31').Value = '  +1.72%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = '0.115'
$ws.Range('E32').Value = '  +0.26%  '
$ws.Range('E33').Value = '  -2.60%  '
$ws.Range('E34').Value = '  +0.16%  '
$ws.Range('D35').Value = '0.0487'
$ws.Range('E35').Value = '  +0.79%  '
$ws.Range('D36').Value = '51.56'
$ws.Range('E36').Value = '  -1.07%  '
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('D38').Value = '3.02'
$ws.Range('E38').Value = '  +1.44%  '
$ws.Range('D39').Value = '3.38'
$ws.Range('E39').Value = '  -0.50%  '
$ws.Range('D40').Value = '0.317'
$ws.Range('E40').Value = '  +11.35%  '
$ws.Range('D41').Value = '140.05'
$ws.Range('E41').Value = '  +3.99%  '
$ws.Range('B42').Value = 'Stellar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D42').Value = '0.125'
$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').Value = '1.98'
$ws.Range('E43').Value = '  +0.96%  '
$ws.Range('E44').Value = '  +1.67%  '
$ws.Range('D45').Value = '16.75'
$ws.Range('E45').Value = '  -0.90%  '
$ws.Range('D46').Value = '2.22'
$ws.Range('E46').Value = '  -0.86%  '
$ws.Range('D47').Value = '21.34'
$ws.Range('E47').Value = '  -1.48%  '
$ws.Range('D48').Value = '2.106.62'
$ws.Range('E48').Value = '  -1.86%  '
$ws.Range('E49').Value = '  -0.86%  '
$ws.Range('D50').Value = '1.92'
$ws.Range('E50').Value = '  +2.27%  '
$ws.Range('D51').Value = '1.73'
$ws.Range('E51').Value = '  +18.05%  '
